$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new dSF (column F) value
$updates = @{
    2  = -5
    3  = 2
    4  = 0
    8  = -1
    9  = 0
    15 = -4
    18 = 3
    22 = -1
    24 = 3
    26 = 1
    27 = -6
    32 = -3
    35 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
